$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 112201331
$ws.Range("B14").Value = 90658
$ws.Range("C14").Value = 'Ovaliderad'
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 4361
$ws.Range("F14").Value = 'Orange taggsvamp'
$ws.Range("G14").Value = 'Hydnellum aurantiacum'
$ws.Range("H14").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("P14").Value = 'Skogalund (Skogalund), Nrk'
$ws.Range("Q14").Value = 531944.4384348277
$ws.Range("R14").Value = 6554005.230760631
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = 'Örebro'
$ws.Range("U14").Value = 'Örebro'
$ws.Range("V14").Value = 'Närke'
$ws.Range("W14").Value = 'Asker'
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = '2023-09-19'
$ws.Range("Z14").Value = '15:46'
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = '2023-09-19'
$ws.Range("AB14").Value = '15:46'
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AW14").Value = 'Magnus Friberg'
$ws.Range("AX14").Value = 'Magnus Friberg'

# Row 15
$ws.Range("A15").Value = 112201510
$ws.Range("B15").Value = 88899
$ws.Range("C15").Value = 'Ovaliderad'
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 3286
$ws.Range("F15").Value = 'Flattoppad klubbsvamp'
$ws.Range("G15").Value = 'Clavariadelphus truncatus'
$ws.Range("H15").Value = '(Quél.) Donk'
$ws.Range("P15").Value = 'Skogalund (Skogalund), Nrk'
$ws.Range("Q15").Value = 531944.4384348277
$ws.Range("R15").Value = 6554005.230760631
$ws.Range("S15").Value = 5
$ws.Range("T15").Value = 'Örebro'
$ws.Range("U15").Value = 'Örebro'
$ws.Range("V15").Value = 'Närke'
$ws.Range("W15").Value = 'Asker'
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = '2023-09-19'
$ws.Range("Z15").Value = '15:52'
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = '2023-09-19'
$ws.Range("AB15").Value = '15:52'
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AG15").Value = $false
$ws.Range("AW15").Value = 'Magnus Friberg'
$ws.Range("AX15").Value = 'Magnus Friberg'

# Row 16
$ws.Range("A16").Value = 112202051
$ws.Range("B16").Value = 90660
$ws.Range("C16").Value = 'Ovaliderad'
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 4362
$ws.Range("F16").Value = 'Blå taggsvamp'
$ws.Range("G16").Value = 'Hydnellum caeruleum'
$ws.Range("H16").Value = '(Hornem.) P.Karst.'
$ws.Range("P16").Value = 'Skogalund (Skogalund), Nrk'
$ws.Range("Q16").Value = 531944.4384348277
$ws.Range("R16").Value = 6554005.230760631
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = 'Örebro'
$ws.Range("U16").Value = 'Örebro'
$ws.Range("V16").Value = 'Närke'
$ws.Range("W16").Value = 'Asker'
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = '2023-09-19'
$ws.Range("Z16").Value = '16:02'
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = '2023-09-19'
$ws.Range("AB16").Value = '16:02'
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AW16").Value = 'Magnus Friberg'
$ws.Range("AX16").Value = 'Magnus Friberg'

# Row 17
$ws.Range("A17").Value = 112202353
$ws.Range("B17").Value = 90659
$ws.Range("C17").Value = 'Ovaliderad'
$ws.Range("D17").Value = 'VU'
$ws.Range("E17").Value = 786
$ws.Range("F17").Value = 'Brandtaggsvamp'
$ws.Range("G17").Value = 'Hydnellum auratile'
$ws.Range("H17").Value = '(Britzelm.) Maas Geest.'
$ws.Range("P17").Value = 'Skogalund (Skogalund), Nrk'
$ws.Range("Q17").Value = 532058.0605805111
$ws.Range("R17").Value = 6553947.659028449
$ws.Range("S17").Value = 5
$ws.Range("T17").Value = 'Örebro'
$ws.Range("U17").Value = 'Örebro'
$ws.Range("V17").Value = 'Närke'
$ws.Range("W17").Value = 'Asker'
$ws.Range("Y17").NumberFormat = "@"
$ws.Range("Y17").Value = '2023-09-19'
$ws.Range("Z17").Value = '16:52'
$ws.Range("AA17").NumberFormat = "@"
$ws.Range("AA17").Value = '2023-09-19'
$ws.Range("AB17").Value = '16:52'
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AW17").Value = 'Magnus Friberg'
$ws.Range("AX17").Value = 'Magnus Friberg'
